$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the column header names in row 1 (shorter labels)
$ws.Range("A1").Value = " Time"
$ws.Range("B1").Value = " ID"
$ws.Range("C1").Value = " Lat"
$ws.Range("D1").Value = " Long"
$ws.Range("E1").Value = " Alt"
$ws.Range("J1").Value = " AntType"
$ws.Range("K1").Value = " Gain"
$ws.Range("L1").Value = " CenterFreq"
$ws.Range("M1").Value = " BandWith"
$ws.Range("N1").Value = " SNR"
$ws.Range("O1").Value = " x"
$ws.Range("P1").Value = " y"
$ws.Range("Q1").Value = " z"

# Header row used to need extra height for the long rotated labels;
# the shorter labels need less room now.
$ws.Range("A1").RowHeight = 60

# Refresh the sample data rows (2-11) with corrected/regenerated values
$ws.Range("B2").Value = -444
$ws.Range("C2").Value = 47.684197573978061
$ws.Range("D2").Value = 9.3970467716033603
$ws.Range("E2").Value = 7388
$ws.Range("F2").Value = -10
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = -56
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 255
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 8500000
$ws.Range("M2").Value = 11000
$ws.Range("N2").Value = 136
$ws.Range("O2").Value = 59
$ws.Range("P2").Value = 32
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = "The quick brown fox jumps over a lazy dog."
$ws.Range("B3").Value = -740
$ws.Range("C3").Value = 47.695732678231174
$ws.Range("D3").Value = 9.4297245355518928
$ws.Range("E3").Value = 1394
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = -89
$ws.Range("H3").Value = -52
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 71
$ws.Range("L3").Value = 8900000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = 134
$ws.Range("O3").Value = 55
$ws.Range("P3").Value = -20
$ws.Range("Q3").Value = -62
$ws.Range("R3").Value = "Vom Ödipuskomplex maßlos gequält, übt Wilfried zyklisches Jodeln."
$ws.Range("B4").Value = 306
$ws.Range("C4").Value = 47.681331740886108
$ws.Range("D4").Value = 9.4206465625590674
$ws.Range("E4").Value = 4347
$ws.Range("F4").Value = -29
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = -5
$ws.Range("I4").Value = 101
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 70
$ws.Range("L4").Value = 9800000
$ws.Range("M4").Value = 16000
$ws.Range("N4").Value = 134
$ws.Range("O4").Value = 51
$ws.Range("P4").Value = -21
$ws.Range("Q4").Value = 40
$ws.Range("R4").Value = "Falsches Üben von Xylophonmusik quält jeden größeren Zwerg."
$ws.Range("B5").Value = -183
$ws.Range("C5").Value = 47.681572705305427
$ws.Range("D5").Value = 9.4165398955014066
$ws.Range("E5").Value = 3334
$ws.Range("F5").Value = -27
$ws.Range("G5").Value = -72
$ws.Range("H5").Value = -51
$ws.Range("I5").Value = 255
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 40
$ws.Range("L5").Value = 10400000
$ws.Range("M5").Value = 13000
$ws.Range("N5").Value = 126
$ws.Range("O5").Value = -60
$ws.Range("P5").Value = -63
$ws.Range("Q5").Value = -10
$ws.Range("R5").Value = "Vom Ödipuskomplex maßlos gequält, übt Wilfried zyklisches Jodeln."
$ws.Range("B6").Value = -506
$ws.Range("C6").Value = 47.691367571696794
$ws.Range("D6").Value = 9.4329133150378581
$ws.Range("E6").Value = 11854
$ws.Range("F6").Value = 69
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = -85
$ws.Range("I6").Value = 103
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 115
$ws.Range("L6").Value = 8800000
$ws.Range("M6").Value = 14000
$ws.Range("N6").Value = 85
$ws.Range("O6").Value = -58
$ws.Range("P6").Value = -32
$ws.Range("Q6").Value = -40
$ws.Range("R6").Value = "The quick brown fox jumps over a lazy dog."
$ws.Range("B7").Value = -64
$ws.Range("C7").Value = 47.673540236618607
$ws.Range("D7").Value = 9.4130118527567816
$ws.Range("E7").Value = 3039
$ws.Range("F7").Value = -30
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = -39
$ws.Range("I7").Value = 103
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 13
$ws.Range("L7").Value = 9900000
$ws.Range("M7").Value = 10000
$ws.Range("N7").Value = 71
$ws.Range("O7").Value = 53
$ws.Range("P7").Value = -35
$ws.Range("Q7").Value = -49
$ws.Range("R7").Value = "Falsches Üben von Xylophonmusik quält jeden größeren Zwerg."
$ws.Range("B8").Value = 198
$ws.Range("C8").Value = 47.713251513591139
$ws.Range("D8").Value = 9.4311152235341602
$ws.Range("E8").Value = 10963
$ws.Range("F8").Value = -78
$ws.Range("G8").Value = 26
$ws.Range("H8").Value = -44
$ws.Range("I8").Value = 103
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 16
$ws.Range("L8").Value = 9300000
$ws.Range("M8").Value = 18000
$ws.Range("N8").Value = 19
$ws.Range("O8").Value = 42
$ws.Range("P8").Value = -65
$ws.Range("Q8").Value = 54
$ws.Range("R8").Value = "Falsches Üben von Xylophonmusik quält jeden größeren Zwerg."
$ws.Range("B9").Value = 261
$ws.Range("C9").Value = 47.716840315494593
$ws.Range("D9").Value = 9.4340819395363518
$ws.Range("E9").Value = 5299
$ws.Range("F9").Value = 57
$ws.Range("G9").Value = 65
$ws.Range("H9").Value = 7
$ws.Range("I9").Value = 102
$ws.Range("J9").Value = 255
$ws.Range("K9").Value = 90
$ws.Range("L9").Value = 8900000
$ws.Range("M9").Value = 19000
$ws.Range("N9").Value = 86
$ws.Range("O9").Value = 36
$ws.Range("P9").Value = 18
$ws.Range("Q9").Value = -1
$ws.Range("R9").Value = "Franz jagt im komplett verwahrlosten Taxi quer durch Bayern."
$ws.Range("B10").Value = 26
$ws.Range("C10").Value = 47.68771994901946
$ws.Range("D10").Value = 9.4059525086396896
$ws.Range("E10").Value = 5903
$ws.Range("F10").Value = 24
$ws.Range("G10").Value = -21
$ws.Range("H10").Value = 60
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 255
$ws.Range("K10").Value = 82
$ws.Range("L10").Value = 8600000
$ws.Range("M10").Value = 14000
$ws.Range("N10").Value = 52
$ws.Range("O10").Value = 33
$ws.Range("P10").Value = -35
$ws.Range("Q10").Value = 30
$ws.Range("R10").Value = "The quick brown fox jumps over a lazy dog."
$ws.Range("B11").Value = -533
$ws.Range("C11").Value = 47.683094400362137
$ws.Range("D11").Value = 9.3999438235127855
$ws.Range("E11").Value = 4600
$ws.Range("F11").Value = 18
$ws.Range("G11").Value = 84
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 255
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 9800000
$ws.Range("M11").Value = 19000
$ws.Range("N11").Value = 71
$ws.Range("O11").Value = -16
$ws.Range("P11").Value = -8
$ws.Range("Q11").Value = -52
$ws.Range("R11").Value = "Franz jagt im komplett verwahrlosten Taxi quer durch Bayern."
